$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E1 is no longer part of the used data (column E drops out of row 1 entirely)
$ws.Range("E1").Clear() | Out-Null

# Clear old text values from column B / A cells that previously held strings
$ws.Range("B2").Value = $null
$ws.Range("B3").Value = $null
$ws.Range("B5").Value = $null
$ws.Range("A6").Value = $null
$ws.Range("B7").Value = $null
$ws.Range("B8").Value = $null
$ws.Range("A10").Value = $null
$ws.Range("B12").Value = $null
$ws.Range("A13").Value = $null

# Set the new row 1 header values (order matches shared-string insertion order)
$ws.Range("A1").Value = "style"
$ws.Range("G1").Value = "font=Aptos"
$ws.Range("B1").Value = "Name=Text"
$ws.Range("C1").Value = "PARENT  =null"
$ws.Range("D1").Value = "coloR=   000000"
$ws.Range("F1").Value = "siZE = 12"

# Copy the style from B1 (s="1") to the new header cells so they share formatting
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1:D1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null

# Update the selection to match the new active cell
$ws.Range("F1").Select() | Out-Null
